$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Crypto price/volume list update (GitHub Actions scheduled refresh)
# D-column price strings are prefixed with a literal apostrophe so Excel
# stores them as text (matching the source data's dotted-thousands format)
# instead of coercing them to numbers.
$ws.Range("D2").Value = "'68.471.43"
$ws.Range("E2").Value = "  -0.88%  "
$ws.Range("D3").Value = "'3.852.83"
$ws.Range("E3").Value = "  -1.89%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "'598.91"
$ws.Range("E5").Value = "  -0.97%  "
$ws.Range("D6").Value = "'167.93"
$ws.Range("E6").Value = "  +1.73%  "
$ws.Range("D7").Value = "'3.852.02"
$ws.Range("E7").Value = "  -1.86%  "
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("D9").Value = "'0.529"
$ws.Range("E9").Value = "  -0.56%  "
$ws.Range("E10").Value = "  -0.92%  "
$ws.Range("D11").Value = "'6.34"
$ws.Range("E11").Value = "  -1.11%  "
$ws.Range("D12").Value = "'0.463"
$ws.Range("E12").Value = "  +0.05%  "
$ws.Range("D13").Value = "'0.0000250"
$ws.Range("E13").Value = "  +1.56%  "
$ws.Range("D14").Value = "'37.42"
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "'4.500.84"
$ws.Range("E15").Value = "  -1.70%  "
$ws.Range("D16").Value = "'3.856.73"
$ws.Range("E16").Value = "  -2.19%  "
$ws.Range("D17").Value = "'68.646.58"
$ws.Range("E17").Value = "  -0.75%  "
$ws.Range("D18").Value = "'7.57"
$ws.Range("E18").Value = "  +1.06%  "
$ws.Range("D19").Value = "'18.22"
$ws.Range("E19").Value = "  +6.00%  "
$ws.Range("E20").Value = "  -1.20%  "
$ws.Range("D21").Value = "'10.75"
$ws.Range("E21").Value = "  -4.01%  "
$ws.Range("D22").Value = "'472.95"
$ws.Range("E22").Value = "  -3.40%  "
$ws.Range("D23").Value = "'0.737"
$ws.Range("E23").Value = "  +1.49%  "
$ws.Range("D24").Value = "'0.0000160"
$ws.Range("E24").Value = "  -4.04%  "
$ws.Range("D25").Value = "'84.60"
$ws.Range("E25").Value = "  +0.05%  "
$ws.Range("D26").Value = "'2.24"
$ws.Range("E26").Value = "  -1.05%  "
$ws.Range("D27").Value = "'12.34"
$ws.Range("E27").Value = "  +1.53%  "
$ws.Range("D28").Value = "'10.06"
$ws.Range("E28").Value = "  -0.85%  "
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("D30").Value = "'2.94"
$ws.Range("E30").Value = "  -0.21%  "
$ws.Range("D31").Value = "'4.003.74"
$ws.Range("E31").Value = "  -1.72%  "
$ws.Range("D32").Value = "'7.79"
$ws.Range("E32").Value = "  -1.37%  "
$ws.Range("D33").Value = "'2.32"
$ws.Range("E33").Value = "  -3.14%  "
$ws.Range("E34").Value = "  -3.82%  "
$ws.Range("D35").Value = "'3.823.87"
$ws.Range("E35").Value = "  -1.09%  "
$ws.Range("E36").Value = "  -1.32%  "
$ws.Range("D37").Value = "'0.140"
$ws.Range("E37").Value = "  +0.42%  "
$ws.Range("D38").Value = "'5.97"
$ws.Range("E38").Value = "  +0.43%  "
$ws.Range("B39").Value = "Mantle"
$ws.Range("C39").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D39").Value = "'1.01"
$ws.Range("E39").Value = "  -2.96%  "
$ws.Range("B40").Value = "dogwifhat"
$ws.Range("C40").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D40").Value = "'3.30"
$ws.Range("E40").Value = "  +8.74%  "
$ws.Range("E41").Value = "  +0.04%  "
$ws.Range("D42").Value = "'0.316"
$ws.Range("E42").Value = "  -1.81%  "
$ws.Range("B43").Value = "Bittensor"
$ws.Range("C43").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D43").Value = "'431.05"
$ws.Range("E43").Value = "  -2.41%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'2.01"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").Value = "'47.58"
$ws.Range("E45").Value = "  -1.88%  "
$ws.Range("D47").Value = "'8.62"
$ws.Range("E47").Value = "  +1.82%  "
$ws.Range("D48").Value = "'0.000274"
$ws.Range("E48").Value = "  +15.00%  "
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D49").Value = "'142.87"
$ws.Range("E49").Value = "  +0.59%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").Value = "'0.0361"
$ws.Range("E50").Value = "  +0.76%  "
$ws.Range("D51").Value = "'39.19"
$ws.Range("E51").Value = "  -0.54%  "

Write-Host "Updated cryptos worksheet"
